$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311
$ws.Range("D22").Value = 1497.908
$ws.Range("D23").Value = 1497.908
